$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily sales rows appended after the existing data (rows 14-36),
# continuing the consecutive-date series and Dollars figures.
$newRows = @(
    @(45566, 220),
    @(45567, 250),
    @(45568, 100),
    @(45569, 120),
    @(45570, "None"),
    @(45571, 500),
    @(45572, 300),
    @(45573, 420),
    @(45574, 380),
    @(45575, 390),
    @(45576, 240),
    @(45577, 220),
    @(45578, 310),
    @(45579, 330),
    @(45580, 360),
    @(45581, 355),
    @(45582, 347),
    @(45583, "None"),
    @(45584, 175),
    @(45585, 185),
    @(45586, 190),
    @(45587, 410),
    @(45588, 405)
)

$r = 14
foreach ($row in $newRows) {
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Trailing rows that only carry the date-column formatting, no values
# (mirrors the format being dragged further down than the data).
for ($r = 37; $r -le 61; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy"
}

$ws.Range("E32").Select() | Out-Null
